# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") was previously populated with a
# "Strike#" style count; this re-derives/writes the correct K values for
# every data row (rows 2-70) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-G ("K") values, in row order starting at row 2.
$kValues = @(
    0,1,1,2,2,1,1,2,1,1,
    1,1,1,0,1,2,2,2,0,2,
    1,1,0,1,2,0,1,2,0,2,
    0,0,0,1,4,0,1,1,0,1,
    2,1,1,1,3,1,2,1,1,1,
    0,1,2,3,0,1,0,0,0,1,
    0,0,0,0,1,1,2,1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
